$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F3").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F4").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("E5").Value = "7,50"
$ws.Range("F5").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F6").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F7").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F8").Value = "11 de jun., 13:51 UTC ·"
$ws.Range("F9").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F10").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F11").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F12").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F13").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F15").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F16").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F17").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F18").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F19").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F20").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F21").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F22").Value = "11 de jun., 13:52 UTC ·"
$ws.Range("F23").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F24").Value = "11 de jun., 13:53 UTC ·"
$ws.Range("F25").Value = "11 de jun., 13:53 UTC ·"
